# Refatoracao - calculos de apoio medio
# Adds std/min/max breakdown columns for "arrecadado" (renamed from *_sucesso),
# "apoio" and "contribuicoes", matching the new panorama.xlsx layout (A1:U4).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): rename/add headers for the new apoio_* and contribuicoes_* breakdown columns ---
$ws.Range("A1").Value = "modalidade"
$ws.Range("B1").Value = "total"
$ws.Range("C1").Value = "total_sucesso"
$ws.Range("D1").Value = "particip"
$ws.Range("E1").Value = "taxa_sucesso"
$ws.Range("F1").Value = "arrecadado_sucesso"
$ws.Range("G1").Value = "arrecadado_avg"
$ws.Range("H1").Value = "arrecadado_std"
$ws.Range("I1").Value = "arrecadado_min"
$ws.Range("J1").Value = "arrecadado_max"
$ws.Range("K1").Value = "apoio_medio"
$ws.Range("L1").Value = "apoio_std"
$ws.Range("M1").Value = "apoio_min"
$ws.Range("N1").Value = "apoio_max"
$ws.Range("O1").Value = "contribuicoes"
$ws.Range("P1").Value = "contribuicoes_med"
$ws.Range("Q1").Value = "contribuicoes_std"
$ws.Range("R1").Value = "contribuicoes_min"
$ws.Range("S1").Value = "contribuicoes_max"
$ws.Range("T1").Value = "menor_ano"
$ws.Range("U1").Value = "maior_ano"

# --- Data rows: shift old columns to their new positions and fill in new statistics ---
# Row 2 (aon)
$ws.Range("A2").Value = "aon"
$ws.Range("B2").Value = 1335
$ws.Range("C2").Value = 830
$ws.Range("D2").Value = 0.3828505878979065
$ws.Range("E2").Value = 0.6217228464419475
$ws.Range("F2").Value = 24063279.82732303
$ws.Range("G2").Value = 28991.90340641329
$ws.Range("H2").Value = 44961.93536949201
$ws.Range("I2").Value = 41.81688448509265
$ws.Range("J2").Value = 679297.6600721752
$ws.Range("K2").Value = 91.85574933975617
$ws.Range("L2").Value = 49.08980856017526
$ws.Range("M2").Value = 13.93896149503088
$ws.Range("N2").Value = 792.0360759681182
$ws.Range("O2").Value = 263553
$ws.Range("P2").Value = 317.533734939759
$ws.Range("Q2").Value = 423.019225146675
$ws.Range("R2").Value = 1
$ws.Range("S2").Value = 6494
$ws.Range("T2").Value = 2011
$ws.Range("U2").Value = 2023

# Row 3 (flex)
$ws.Range("A3").Value = "flex"
$ws.Range("B3").Value = 1468
$ws.Range("C3").Value = 1383
$ws.Range("D3").Value = 0.4209922569544021
$ws.Range("E3").Value = 0.9420980926430518
$ws.Range("F3").Value = 18362131.9375591
$ws.Range("G3").Value = 13277.02960054888
$ws.Range("H3").Value = 33934.82811955066
$ws.Range("I3").Value = 10.77163914429046
$ws.Range("J3").Value = 708972.7845446636
$ws.Range("K3").Value = 77.41063997458096
$ws.Range("L3").Value = 39.50983355883143
$ws.Range("M3").Value = 10.77163914429046
$ws.Range("N3").Value = 461.5197709071476
$ws.Range("O3").Value = 203646
$ws.Range("P3").Value = 147.2494577006508
$ws.Range("Q3").Value = 327.6748910926806
$ws.Range("R3").Value = 1
$ws.Range("S3").Value = 7954
$ws.Range("T3").Value = 2016
$ws.Range("U3").Value = 2023

# Row 4 (sub)
$ws.Range("A4").Value = "sub"
$ws.Range("B4").Value = 684
$ws.Range("C4").Value = 152
$ws.Range("D4").Value = 0.1961571551476914
$ws.Range("E4").Value = 0.2222222222222222
$ws.Range("F4").Value = 43186.9577547848
$ws.Range("G4").Value = 284.1247220709527
$ws.Range("H4").Value = 650.5808076401024
$ws.Range("I4").Value = 1.087396962410123
$ws.Range("J4").Value = 5087.076865717208
$ws.Range("K4").Value = 21.28348419490777
$ws.Range("L4").Value = 15.01968006252796
$ws.Range("M4").Value = 1.011042153300025
$ws.Range("N4").Value = 84.0771316599004
$ws.Range("O4").Value = 2208
$ws.Range("P4").Value = 14.52631578947368
$ws.Range("Q4").Value = 31.86830254134198
$ws.Range("R4").Value = 1
$ws.Range("S4").Value = 208
$ws.Range("T4").Value = 2016
$ws.Range("U4").Value = 2023

# --- Number formats for data rows (rows 2:4) matching column semantics ---
$ws.Range("B2:C4").NumberFormat = "#,##0"
$ws.Range("D2:E4").NumberFormat = "0.00%"
$ws.Range("F2:N4").NumberFormat = "R$ #,##0.00"
$ws.Range("O2:S4").NumberFormat = "#,##0"
